# Update the "All" labels in column A to "Combined" for the race-strata
# summary rows of the Antonio table (every 3rd row starting at row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 5, 8, 11, 14, 17)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "Combined"
}
